# Adding element/charge balance validation test data:
# update empirical formulas and charges on the "Species types" sheet so
# the fixture exercises the new element/charge balance validator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Species types")

# Empirical formula (column D)
$ws.Range("D2").Value = "NaHCO"
$ws.Range("D3").Value = "N2O2P"
$ws.Range("D4").Value = "N4O4P2"
$ws.Range("D5").Value = "N10O10P5"
$ws.Range("D6").Value = "N5O5"
$ws.Range("D7").Value = "P5"

# Charge (column F)
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -5
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = 1

# Make "Species types" the active sheet (matches workbookView activeTab=4)
$ws.Activate()
